# BOM.xlsx minor update
# 1. Fix the Ref Des ordering for the D-series parts (row 4 / item 3)
# 2. Correct the Description/Value for the same row to the LED description
# 3. Remove the trailing SUM total row (row 9) that is no longer needed
# 4. Leave selection on F5 as last edited by the author

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "D1,D2,D3,D4,D5,D6,D7,D8,D9,D10,D11"
$ws.Range("F4").Value = "LED GREEN CLEAR 0603 SMD"

$ws.Range("C9").EntireRow.Delete()

$ws.Range("F5").Select()
